# Update "Pais" COVID dashboard sheet to the 08:22 snapshot:
#  - bump the "updated at" timestamp (07:52 -> 08:22)
#  - refresh case numbers for the countries whose ranking moved since
#    the 07:52 snapshot (their row now holds a different country, with
#    that country's updated stats)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header (row 1, col A) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 7 de Abril de 2020 a las 08:22"

# --- Refreshed rows (Pais, Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 17 -> Austria
$ws.Cells.Item(17, 1).Value = "Austria"
$ws.Cells.Item(17, 2).Value = 12332
$ws.Cells.Item(17, 3).Value = 35
$ws.Cells.Item(17, 5).Value = 8649

# Row 85 -> Uzbekistan (was Costa Rica)
$ws.Cells.Item(85, 1).Value = "Uzbekistan"
$ws.Cells.Item(85, 2).Value = 472
$ws.Cells.Item(85, 3).Value = 15
$ws.Cells.Item(85, 4).Value = 30
$ws.Cells.Item(85, 5).Value = 440
$ws.Cells.Item(85, 6).Value = 8

# Row 86 -> Costa Rica (was Republica de Chipre)
$ws.Cells.Item(86, 1).Value = "Costa Rica"
$ws.Cells.Item(86, 2).Value = 467
$ws.Cells.Item(86, 4).Value = 18
$ws.Cells.Item(86, 5).Value = 447
$ws.Cells.Item(86, 6).Value = 14
$ws.Cells.Item(86, 8).Value = 2

# Row 87 -> Republica de Chipre (was Uzbekistan)
$ws.Cells.Item(87, 1).Value = "Republica de Chipre"
$ws.Cells.Item(87, 2).Value = 465
$ws.Cells.Item(87, 4).Value = 45
$ws.Cells.Item(87, 5).Value = 411
$ws.Cells.Item(87, 6).Value = 11
$ws.Cells.Item(87, 8).Value = 9

# Row 91 -> Taiwan
$ws.Cells.Item(91, 1).Value = "Taiwan"
$ws.Cells.Item(91, 2).Value = 376
$ws.Cells.Item(91, 3).Value = 3
$ws.Cells.Item(91, 4).Value = 61
$ws.Cells.Item(91, 5).Value = 310

# Row 92 -> Oman (was Burkina Faso)
$ws.Cells.Item(92, 1).Value = "Oman"
$ws.Cells.Item(92, 2).Value = 371
$ws.Cells.Item(92, 3).Value = 40
$ws.Cells.Item(92, 4).Value = 67
$ws.Cells.Item(92, 5).Value = 302
$ws.Cells.Item(92, 6).Value = 3
$ws.Cells.Item(92, 8).Value = 2

# Row 93 -> Burkina Faso (was Cuba)
$ws.Cells.Item(93, 1).Value = "Burkina Faso"
$ws.Cells.Item(93, 2).Value = 364
$ws.Cells.Item(93, 3).Value = 0
$ws.Cells.Item(93, 4).Value = 108
$ws.Cells.Item(93, 5).Value = 238
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 8).Value = 18

# Row 94 -> Cuba (was Reunion)
$ws.Cells.Item(94, 1).Value = "Cuba"
$ws.Cells.Item(94, 2).Value = 363
$ws.Cells.Item(94, 3).Value = 13
$ws.Cells.Item(94, 4).Value = 18
$ws.Cells.Item(94, 5).Value = 336
$ws.Cells.Item(94, 6).Value = 12
$ws.Cells.Item(94, 8).Value = 9

# Row 95 -> Reunion (was Jordania)
$ws.Cells.Item(95, 1).Value = "Reunion"
$ws.Cells.Item(95, 4).Value = 40
$ws.Cells.Item(95, 5).Value = 309
$ws.Cells.Item(95, 6).Value = 4
$ws.Cells.Item(95, 8).Value = 0

# Row 96 -> Jordania (was Oman)
$ws.Cells.Item(96, 1).Value = "Jordania"
$ws.Cells.Item(96, 2).Value = 349
$ws.Cells.Item(96, 4).Value = 126
$ws.Cells.Item(96, 5).Value = 217
$ws.Cells.Item(96, 6).Value = 5
$ws.Cells.Item(96, 8).Value = 6

# Row 110 -> Georgia (was Bolivia)
$ws.Cells.Item(110, 1).Value = "Georgia"
$ws.Cells.Item(110, 2).Value = 195
$ws.Cells.Item(110, 3).Value = 7
$ws.Cells.Item(110, 4).Value = 39
$ws.Cells.Item(110, 5).Value = 154
$ws.Cells.Item(110, 6).Value = 6
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = 2

# Row 111 -> Bolivia (was Georgia)
$ws.Cells.Item(111, 1).Value = "Bolivia"
$ws.Cells.Item(111, 2).Value = 194
$ws.Cells.Item(111, 3).Value = 11
$ws.Cells.Item(111, 4).Value = 2
$ws.Cells.Item(111, 5).Value = 178
$ws.Cells.Item(111, 6).Value = 3
$ws.Cells.Item(111, 7).Value = 3
$ws.Cells.Item(111, 8).Value = 14

# Row 113 -> Sri Lanka
$ws.Cells.Item(113, 1).Value = "Sri Lanka"
$ws.Cells.Item(113, 2).Value = 180
$ws.Cells.Item(113, 3).Value = 2
$ws.Cells.Item(113, 5).Value = 136
